$d = $word.ActiveDocument

# 1) Move the "_GoBack" bookmark from its old location (right after the
#    word "Recibe") to a new location in the middle of the run
#    "modalidad_pago_reserva", splitting it into "modalidad_pago_" and
#    "reserva". Locate the text precisely via Find so we do not depend on
#    hard-coded character offsets.
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found1 = $findRange.Find.Execute("modalidad_pago_reserva", $false, $false, $false,
                                   $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate 'modalidad_pago_reserva' in the document."
}
$prefixText = "modalidad_pago_"
$splitPoint = $findRange.Start + $prefixText.Length
$bookmarkRange = $d.Range($splitPoint, $splitPoint)
[void]$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# 2) Remove the trailing block of paragraphs that followed the lone "RNE"
#    paragraph (from "-Cuando se confirma el pago..." through "...es de
#    rango "Oficinista" o "Gerente"." at the very end of the document body).
$startText = "-Cuando se confirma el pago"
$endText = "-El Empleado que Escribe la solicitud de reserva es de rango " + [char]34 + "Oficinista" + [char]34 + " o " + [char]34 + "Gerente" + [char]34 + "."

$startRange = $d.Content
$startRange.Find.ClearFormatting()
$found2 = $startRange.Find.Execute($startText, $false, $false, $false, $false, $false,
                                    $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the start of the block to remove."
}

$endRange = $d.Content
$endRange.Find.ClearFormatting()
$found3 = $endRange.Find.Execute($endText, $false, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not locate the end of the block to remove."
}

# +1 so the final paragraph mark (not included in the matched text) is
# swept up too - this removes the paragraphs entirely instead of leaving
# empty shells behind.
$deleteRange = $d.Range($startRange.Start, $endRange.End + 1)
$deleteRange.Delete()
